# Updates the cryptos list (rows 2-51) with the latest scraped values.
# Column D ("Price") is forced to Text format first because several values
# (e.g. "29.315.58", "1.002") look numeric but must be kept as the literal
# strings from the source feed (they are not valid financial numbers -
# some even use "." as a thousands separator).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.315.58"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").Value = "1.860.69"
$ws.Range("E3").Value = "  +0.06%  "

# Row 4
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "0.7019"
$ws.Range("E5").Value = "  -0.11%  "

# Row 6
$ws.Range("D6").Value = "238.15"
$ws.Range("E6").Value = "  +0.31%  "

# Row 7
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("D8").Value = "0.07868"
$ws.Range("E8").Value = "  +1.51%  "

# Row 9
$ws.Range("D9").Value = "0.3047"
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").Value = "24.81"
$ws.Range("E10").Value = "  +6.57%  "

# Row 11
$ws.Range("D11").Value = "0.08179"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12
$ws.Range("D12").Value = "1.873.35"
$ws.Range("E12").Value = "  +1.15%  "

# Row 13
$ws.Range("D13").Value = "5.219"
$ws.Range("E13").Value = "  +0.96%  "

# Row 14
$ws.Range("D14").Value = "0.7138"
$ws.Range("E14").Value = "  -0.53%  "

# Row 15
$ws.Range("D15").Value = "89.47"
$ws.Range("E15").Value = "  +0.35%  "

# Row 16
$ws.Range("D16").Value = "29.542.57"
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("D17").Value = "5.826"
$ws.Range("E17").Value = "  +0.92%  "

# Row 18
$ws.Range("D18").Value = "0.000007796"
$ws.Range("E18").Value = "  +0.79%  "

# Row 19
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "13.23"
$ws.Range("E19").Value = "  -1.07%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "238.75"
$ws.Range("E20").Value = "  +0.97%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.205.85"
$ws.Range("E21").Value = "  +4.94%  "

# Row 22
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.29%  "

# Row 23
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.20%  "

# Row 24
$ws.Range("D24").Value = "7.556"
$ws.Range("E24").Value = "  +1.46%  "

# Row 25
$ws.Range("D25").Value = "162.93"
$ws.Range("E25").Value = "  +0.79%  "

# Row 26
$ws.Range("D26").Value = "8.899"
$ws.Range("E26").Value = "  -1.09%  "

# Row 27
$ws.Range("D27").Value = "0.1422"
$ws.Range("E27").Value = "  -3.42%  "

# Row 28
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  +0.33%  "

# Row 29
$ws.Range("D29").Value = "1.898"
$ws.Range("E29").Value = "  -5.66%  "

# Row 30
$ws.Range("D30").Value = "1.373"
$ws.Range("E30").Value = "  -4.31%  "

# Row 31
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").Value = "4.312"
$ws.Range("E32").Value = "  -2.40%  "

# Row 33
$ws.Range("D33").Value = "4.049"
$ws.Range("E33").Value = "  +0.12%  "

# Row 34
$ws.Range("D34").Value = "0.05174"
$ws.Range("E34").Value = "  -0.91%  "

# Row 35
$ws.Range("D35").Value = "1.177"

# Row 36
$ws.Range("D36").Value = "0.7075"
$ws.Range("E36").Value = "  +0.10%  "

# Row 37
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  +0.15%  "

# Row 38
$ws.Range("D38").Value = "2.678"
$ws.Range("E38").Value = "  +0.30%  "

# Row 39
$ws.Range("D39").Value = "0.01844"
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("D40").Value = "2.694"
$ws.Range("E40").Value = "  -0.97%  "

# Row 41
$ws.Range("D41").Value = "1.170.01"
$ws.Range("E41").Value = "  +2.61%  "

# Row 42
$ws.Range("D42").Value = "0.9210"
$ws.Range("E42").Value = "  -1.06%  "

# Row 43
$ws.Range("D43").Value = "6.037"
$ws.Range("E43").Value = "  +1.90%  "

# Row 44
$ws.Range("D44").Value = "71.64"
$ws.Range("E44").Value = "  +1.22%  "

# Row 45
$ws.Range("D45").Value = "0.4259"
$ws.Range("E45").Value = "  -0.38%  "

# Row 46
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.12%  "

# Row 47
$ws.Range("D47").Value = "101.88"
$ws.Range("E47").Value = "  -1.63%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5350"
$ws.Range("E48").Value = "  -1.77%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.753"
$ws.Range("E49").Value = "  -2.25%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.157"
$ws.Range("E50").Value = "  -0.19%  "

# Row 51
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "6.977"
$ws.Range("E51").Value = "  -0.02%  "
